# correct the error time table in document
#
# The schedule table has a "第一期评审" (phase-1 review) row running
# 2013.4.4 - 2013.4.4, and the following "第二期开发" (phase-2 dev) row
# starting 2013.4.4. These are off by one: the review should end (and the
# next phase should start) on 2013.4.3, not 2013.4.4.
#
# wdReplaceOne / wdFindContinue / wdCollapseEnd style constants used below:
#   Find.Execute(..., Replace:=1) -> wdReplaceOne (replace just the found hit)
#   Find.Execute(..., Wrap:=1)    -> wdFindContinue

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "2013.4.4-" (start of the 第一期评审 date range, first paragraph of
#    the cell) -> "2013.4.3-"
# ---------------------------------------------------------------------
$rng1 = $d.Range(0, $d.Content.End)
$rng1.Find.Execute("2013.4.4-", $true, $false, $false, $false, $false, $true, 1, $false, "2013.4.3-", 1) | Out-Null

# ---------------------------------------------------------------------
# 2) "2013.4.4" (end of the 第一期评审 date range, second paragraph of the
#    same cell, written as the runs "2013.4" + "." + "4") -> change just
#    the trailing "4" to "3", leaving the rest of the text untouched.
# ---------------------------------------------------------------------
$rng2 = $d.Range(0, $d.Content.End)
$rng2.Find.Execute("2013.4.4", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$lastDigit = $d.Range($rng2.End - 1, $rng2.End)
$lastDigit.Text = "3"

# ---------------------------------------------------------------------
# 3) "2013.4.4" (start of the 第二期开发 date range, immediately followed
#    by a separate "-" run) -> "2013.4.3". The document's "_GoBack" last-
#    edit bookmark is relocated here, right after the corrected text,
#    matching where the author's cursor last made an edit.
# ---------------------------------------------------------------------
$rng3 = $d.Range(0, $d.Content.End)
$rng3.Find.Execute("2013.4.4", $true, $false, $false, $false, $false, $true, 1, $false, "2013.4.3", 1) | Out-Null

$d.Bookmarks.Add("_GoBack", $d.Range($rng3.End, $rng3.End)) | Out-Null
